$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = "[R] Add an index (or counter) to a dataframe by group"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/R-Add-an-index-or-counter-to-a-dataframe-by-group"

# Row 9
$ws.Range("D9").Value = "Data Science 석사하려고 그전에 석사 학위를 1개 더 해야 한다구요?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/why-two-masters/#utm_source=rss&utm_medium=rss&utm_campaign=why-two-masters"

# Row 29
$ws.Range("D29").Value = "[MRI basics] k-space 변화에 따른  이미지 변화 확인하기"
$ws.Range("E29").Value = "https://blog.promedius.ai/mri-basics-k-space-byeonhwae-ddareun-imiji-byeonhwa-hwaginhagi-2/"

# Row 51
$ws.Range("D51").Value = "[javascript] 문자열의 길이를 알려주는 length 속성"
$ws.Range("E51").Value = "https://bskyvision.com/1153"
